# Add a new bulleted accomplishment line right after the
# "Leveraged Elixir to build webhooks..." bullet (end of the Lyft job's
# bullet list, just before the "Exporta.io" entry).
#
# The original paragraph has a trailing empty run (just <w:rPr><w:rtl/></w:rPr>)
# after its text run.  Using Find/Replace with a literal paragraph-mark
# ("^p") in the replacement text both (a) splits the text into two
# paragraphs and (b) naturally absorbs/clears that trailing empty run,
# matching the target structure, where the new paragraph inherits the
# same list numbering / indentation (numPr ilvl=0 numId=1, ind left=720
# hanging=360) from the paragraph it was split from.

$d = $word.ActiveDocument

$oldText = "Leveraged Elixir to build webhooks that communicate with our customers regarding rides in progress"
$newBullet = "Built a GraphQL microservice in Clojure that communicates over Confluent Kafka"

$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, `
    "$oldText^p$newBullet", 2)

if (-not $found) {
    throw "Could not find the 'Leveraged Elixir...' paragraph to split."
}

# The newly created paragraph inherits sz/szCs=22 via the paragraph's
# pPr/rPr, but the new run itself is created without explicit run-level
# sz/szCs (only rtl). Set the font size explicitly (11pt = sz 22) on the
# new run so it serializes with <w:sz w:val="22"/><w:szCs w:val="22"/>,
# matching the sibling bullet runs in this list.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$newBullet*") {
        $p.Range.Font.Size = 11
        $p.Range.Font.SizeBi = 11
        break
    }
}
